$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.935.05"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "2.287.15"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.04"
$ws.Range("E5").Value = "  -5.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.35"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").Value = "2.284.83"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.329"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.42"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "2.694.34"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "57.873.29"
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "2.264.34"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.50"
$ws.Range("E19").Value = "  -4.30%  "
$ws.Range("E20").Value = "  -5.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.37"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.36"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.97"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("E28").Value = "  -6.68%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  -5.84%  "
$ws.Range("D31").Value = "0.0₃0716"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.71"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("E34").Value = "  -5.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.71"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("E39").Value = "  -4.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.47"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.49"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.08"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "285.99"
$ws.Range("E43").Value = "  -8.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0494"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.91"
$ws.Range("E48").Value = "  -5.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.91"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  -0.59%  "
